# Adds a new "Low Caliber Towers" research entry to the localization sheet:
#   - gui/menu/research/description/towers_lowcaliber
#   - gui/menu/research/name/towers_lowcaliber
#
# The sheet (rebalance_localizations) is sorted alphabetically by column A,
# so the two new rows land in their alphabetically-correct spots:
#   * the "description/..." row goes in right before "description/well_contruction"
#   * the "name/..." row goes in right before "name/well_contruction"
# which, in the original (pre-edit) row numbering, means inserting new rows
# at row 196 and (after that first insert shifts things down) row 237.

# 0) give the new research entry a readable name rather than repeating
#    literal strings everywhere below.
$descKey = "gui/menu/research/description/towers_lowcaliber"
$descVal = "Provides basic defense towers utilizing low caliber ammunition."
$nameKey = "gui/menu/research/name/towers_lowcaliber"
$nameVal = "Low Caliber Towers"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rebalance_localizations")

# 1) Insert the description row just above the old row 196
#    ("gui/menu/research/description/well_contruction"), and the name row
#    just above what is then old row 236 ("gui/menu/research/name/well_contruction"),
#    i.e. new row 237, after the shift caused by the first insertion above.
$ws.Rows.Item(196).Insert()
$ws.Rows.Item(237).Insert()

# 2) Fill in the two new rows column-by-column (name row's key, then
#    description row's key, then name row's value, then description row's
#    value) to mirror how the change was authored.
$ws.Range("A237").Value = $nameKey
$ws.Range("A196").Value = $descKey
$ws.Range("B237").Value = $nameVal
$ws.Range("B196").Value = $descVal

# Match the author's final selection/view state as closely as this runtime
# allows (freeze-pane scroll offset isn't independently controllable here,
# but the active cell selection is).
$ws.Range("C18").Select()
